$d = $word.ActiveDocument

# Replace " a short description of the web app. " with the extended sentence
# mentioning the Inspirational Quotes API.
$d.Content.Find.Execute(
    "short description of the web app. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "short description of the web app, as well as an inspirational quote from the Inspirational Quotes API. ",
    2
)
